# BOM.xlsx update: pcb adapted + BOM completed
# Adds new rows for the WiFi antenna kit (router, antennas, cables, adapters)
# just above the existing "Microwave Coaxial adapter" line, and renames the
# router line item.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Make room: insert 3 blank rows before the old row 40 -----------------
# (old row 40 "Murata / Microwave Coaxial adapter" ends up at row 43, and the
#  rows below it all shift down by 3 accordingly)
$ws.Rows("40:42").Insert()

# Copy the formatting of row 39 (Digitec / Asus router line) onto the three
# freshly inserted blank rows so they match the surrounding table styling.
$ws.Range("A39:H39").Copy()
$ws.Range("A40:H42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2 & 3. Rename the router line and fill the new antenna-kit rows --------
# Text cells are written in the same order the original author typed them in
# (manufacturer/description first, order-codes after) so the shared-string
# table comes out in the same sequence as the real edit.
$ws.Range("A40").Value = "Distrelec"
$ws.Range("D40").Value = "Taoglas"
$ws.Range("E39").Value = "Wifi Router - Asus RT-AX86U Pro"
$ws.Range("C40").Value = "302-20-253"
$ws.Range("E40").Value = "Dual Band wifi antenna"
$ws.Range("C42").Value = "301-31-595"
$ws.Range("D42").Value = "Nedis"
$ws.Range("E42").Value = "Antenna cable SMA 2m"
$ws.Range("C41").Value = "`tADP-SMAM-RPSF-G-ND"
$ws.Range("D41").Value = "Linx Technologies"
$ws.Range("E41").Value = "SMA to RP-SMA adapter"

# Remaining (numeric / formula) cells - order does not affect the shared
# string table.
$ws.Range("B40").Value = 1
$ws.Range("F40").Value = 15.64
$ws.Range("G40").Formula = "=F40*B40"

$ws.Range("A41").Value = "Digikey"
$ws.Range("B41").Value = 1
$ws.Range("F41").Value = 6.72
$ws.Range("G41").Formula = "=F41*B41"

$ws.Range("A42").Value = "Distrelec"
$ws.Range("B42").Value = 1
$ws.Range("F42").Value = 10.29
$ws.Range("G42").Formula = "=F42*B42"

# --- 4. Restore the sheet view scroll position described in the diff --------
$ws.Range("I39").Select()
$excel.ActiveWindow.ScrollRow = 16
